$p = $ppt.ActivePresentation

# Layouts used by the new slides.
$layoutSectionHeader = $p.SlideMaster.CustomLayouts.Item(3)   # "Section Header" (title + body)
$layoutTitleContent  = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# --- Slide 2: "Demo" --------------------------------------------------
$s2 = $p.Slides.AddSlide(2, $layoutSectionHeader)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Demo"
# Body placeholder stays empty (left untouched, same as the authored slide).

# --- Slide 3: "Why?" ----------------------------------------------------
$s3 = $p.Slides.AddSlide(3, $layoutTitleContent)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Why?"
# Content placeholder stays empty.

# --- Slide 4: "Code" -----------------------------------------------------
$s4 = $p.Slides.AddSlide(4, $layoutTitleContent)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Code"
# Content placeholder stays empty.

# --- Slide 5: "Difficult Parts" ------------------------------------------
$s5 = $p.Slides.AddSlide(5, $layoutTitleContent)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Difficult Parts"
$tf5 = $s5.Shapes.Item(2).TextFrame
$tf5.TextRange.Text = "Xamarin itself
Despite the platform’s age, tools seem quite buggy
Incomplete documentation"
$tf5.TextRange.Paragraphs(2,1).IndentLevel = 2
$tf5.TextRange.Paragraphs(3,1).IndentLevel = 2

# --- Slide 6: "Highlights" ------------------------------------------------
$s6 = $p.Slides.AddSlide(6, $layoutTitleContent)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Highlights"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Library for interfacing with the NXT is a reusable component"

# --- Slide 7: "Shortcomings" ----------------------------------------------
$s7 = $p.Slides.AddSlide(7, $layoutTitleContent)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Shortcomings"
$tf7 = $s7.Shapes.Item(2).TextFrame
$tf7.TextRange.Text = "Didn’t have time to add EV3 support
Would have liked to test on more devices to verify cross-compatibility actually works"

# --- Slide 8: "Summary" ----------------------------------------------------
$s8 = $p.Slides.AddSlide(8, $layoutTitleContent)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Summary"
$tf8 = $s8.Shapes.Item(2).TextFrame
$tf8.TextRange.Text = "Source code is on GitHub!
URL
Visual Studio project
LDraw"
$tf8.TextRange.Paragraphs(2,1).IndentLevel = 2
$tf8.TextRange.Paragraphs(3,1).IndentLevel = 2
$tf8.TextRange.Paragraphs(4,1).IndentLevel = 2
$tf8.TextRange.Paragraphs(4,1).InsertAfter(" robot model") | Out-Null

# --- Slide 9: "Questions?" --------------------------------------------
$s9 = $p.Slides.AddSlide(9, $layoutSectionHeader)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"
# Body placeholder stays empty.
